$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 483.23077
$ws.Range("J6").Value = 374.75
$ws.Range("L6").Value = 1124.25
$ws.Range("N6").Value = -1348.25

# Row 17
$ws.Range("H17").Value = 671099.9
$ws.Range("J17").Value = 745555.4399999999
$ws.Range("L17").Value = 2236666.32
$ws.Range("N17").Value = -2237002.32

# Row 32
$ws.Range("H32").Value = 29413550
$ws.Range("I32").Value = 1361.1428
$ws.Range("J32").Value = 50002084
$ws.Range("K32").Value = 1361.1428
$ws.Range("L32").Value = 50002084
$ws.Range("M32").Value = -1035.1428
$ws.Range("N32").Value = -50002736

# Row 53
$ws.Range("H53").Value = 52631916
$ws.Range("I53").Value = 150.55556
$ws.Range("K53").Value = 150.55556
$ws.Range("M53").Value = 486.44444

# Row 74
$ws.Range("H74").Value = 3432.6155
$ws.Range("J74").Value = 1185
$ws.Range("L74").Value = 1185
$ws.Range("N74").Value = -3057

# Row 77
$ws.Range("H77").Value = 3432.6155
$ws.Range("J77").Value = 1185
$ws.Range("L77").Value = 5925
$ws.Range("N77").Value = -15285

# Row 98
$ws.Range("H98").Value = 3322.7368
$ws.Range("J98").Value = 5263.6
$ws.Range("L98").Value = 5263.6
$ws.Range("N98").Value = -8259.6

# Row 106
$ws.Range("H106").Value = 1190.4324
$ws.Range("I106").Value = 1221.8667
$ws.Range("J106").Value = 1055.7142
$ws.Range("K106").Value = 1221.8667
$ws.Range("L106").Value = 1055.7142
$ws.Range("M106").Value = -590.8667
$ws.Range("N106").Value = -2317.7142

# Row 122
$ws.Range("H122").Value = 3322.7368
$ws.Range("J122").Value = 5263.6
$ws.Range("L122").Value = 15790.8
$ws.Range("N122").Value = -20690.8

# Row 138
$ws.Range("H138").Value = 4701.5625
$ws.Range("I138").Value = 10933.9
$ws.Range("J138").Value = 3547.426
$ws.Range("K138").Value = 32801.7
$ws.Range("L138").Value = 10642.278
$ws.Range("M138").Value = -27661.7
$ws.Range("N138").Value = -20922.278


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1088.6666
$ws.Range("I2").Value = 1088.6666
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1088.6666
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -975.6666
$ws.Range("N2").ClearContents()

# Row 17
$ws.Range("H17").Value = 5655
$ws.Range("J17").Value = 5655
$ws.Range("L17").Value = 5655
$ws.Range("N17").Value = -6001

# Row 45
$ws.Range("H45").Value = 39056
$ws.Range("I45").Value = 47396
$ws.Range("J45").Value = 2360
$ws.Range("K45").Value = 47396
$ws.Range("L45").Value = 2360
$ws.Range("M45").Value = -47019
$ws.Range("N45").Value = -3114

# Row 64
$ws.Range("H64").Value = 34000
$ws.Range("J64").Value = 34000
$ws.Range("L64").Value = 34000
$ws.Range("N64").Value = -34496

# Row 67
$ws.Range("H67").Value = 34000
$ws.Range("J67").Value = 34000
$ws.Range("L67").Value = 34000
$ws.Range("N67").Value = -35716

# Row 116
$ws.Range("H116").Value = 1088.6666
$ws.Range("I116").Value = 1088.6666
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1088.6666
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1205.3334
$ws.Range("N116").ClearContents()

# Row 122
$ws.Range("H122").Value = 1023.8333
$ws.Range("I122").Value = 1034.7273
$ws.Range("J122").Value = 904
$ws.Range("K122").Value = 3104.1819
$ws.Range("L122").Value = 2712
$ws.Range("M122").Value = -654.1819
$ws.Range("N122").Value = -7612

# Row 132
$ws.Range("H132").Value = 2623.4211
$ws.Range("I132").Value = 2557
$ws.Range("K132").Value = 7671
$ws.Range("M132").Value = -5141


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1088.6666
$ws.Range("I3").Value = 1088.6666
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1088.6666
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -974.6666
$ws.Range("N3").ClearContents()

# Row 134
$ws.Range("H134").Value = 18751954
$ws.Range("I134").Value = 1645.65
$ws.Range("K134").Value = 4936.950000000001
$ws.Range("M134").Value = -2401.950000000001


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 20410832
$ws.Range("I16").Value = 23812304
$ws.Range("K16").Value = 23812304
$ws.Range("M16").Value = -23812017

# Row 22
$ws.Range("H22").Value = 759.5306399999999
$ws.Range("I22").Value = 738.383
$ws.Range("J22").Value = 1256.5
$ws.Range("K22").Value = 738.383
$ws.Range("L22").Value = 1256.5
$ws.Range("M22").Value = -388.383
$ws.Range("N22").Value = -1956.5

# Row 53
$ws.Range("H53").Value = 53684
$ws.Range("J53").Value = 53684
$ws.Range("L53").Value = 53684
$ws.Range("N53").Value = -54898

# Row 58
$ws.Range("H58").Value = 1857.375
$ws.Range("I58").Value = 2237.125
$ws.Range("J58").Value = 1477.625
$ws.Range("K58").Value = 2237.125
$ws.Range("L58").Value = 1477.625
$ws.Range("M58").Value = -2034.125
$ws.Range("N58").Value = -1883.625

# Row 99
$ws.Range("H99").Value = 4002599
$ws.Range("I99").Value = 6668666.5
$ws.Range("K99").Value = 6668666.5
$ws.Range("M99").Value = -6667168.5

# Row 113
$ws.Range("H113").Value = 20410832
$ws.Range("I113").Value = 23812304
$ws.Range("K113").Value = 23812304
$ws.Range("M113").Value = -23810134

# Row 126
$ws.Range("H126").Value = 4002599
$ws.Range("I126").Value = 6668666.5
$ws.Range("K126").Value = 20005999.5
$ws.Range("M126").Value = -20003529.5

# Row 134
$ws.Range("H134").Value = 1393.0883
$ws.Range("I134").Value = 1309.8572
$ws.Range("K134").Value = 3929.5716
$ws.Range("M134").Value = -1394.5716

# Row 136
$ws.Range("H136").Value = 1857.375
$ws.Range("I136").Value = 2237.125
$ws.Range("J136").Value = 1477.625
$ws.Range("K136").Value = 6711.375
$ws.Range("L136").Value = 4432.875
$ws.Range("M136").Value = -4161.375
$ws.Range("N136").Value = -9532.875


$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 28477.166
$ws.Range("I12").Value = 67.5
$ws.Range("J12").Value = 36594.215
$ws.Range("K12").Value = 202.5
$ws.Range("L12").Value = 109782.645
$ws.Range("M12").Value = -29.5
$ws.Range("N12").Value = -110128.645

# Row 113
$ws.Range("H113").Value = 394.375
$ws.Range("I113").Value = 360.33334
$ws.Range("J113").Value = 402.23077
$ws.Range("K113").Value = 1081.00002
$ws.Range("L113").Value = 1206.69231
$ws.Range("M113").Value = 1088.99998
$ws.Range("N113").Value = -5546.69231

# Row 115
$ws.Range("H115").Value = 1574.5555
$ws.Range("I115").Value = 396
$ws.Range("J115").Value = 3931.6667
$ws.Range("K115").Value = 1188
$ws.Range("L115").Value = 11795.0001
$ws.Range("M115").Value = -13
$ws.Range("N115").Value = -14145.0001

# Row 131
$ws.Range("H131").Value = 52580.65
$ws.Range("I131").Value = 1949.4546
$ws.Range("K131").Value = 5848.3638
$ws.Range("M131").Value = -808.3638000000001


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2255.7917
$ws.Range("I102").Value = 1872.7222
$ws.Range("K102").Value = 1872.7222
$ws.Range("M102").Value = -250.7221999999999


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6834.407
$ws.Range("I7").Value = 2964.7693
$ws.Range("K7").Value = 2964.7693
$ws.Range("M7").Value = -2852.7693

# Row 30
$ws.Range("H30").Value = 7000
$ws.Range("I30").Value = 7000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -6892
$ws.Range("N30").ClearContents()

# Row 46
$ws.Range("H46").Value = 14366.333
$ws.Range("J46").Value = 6924.875
$ws.Range("L46").Value = 6924.875
$ws.Range("N46").Value = -7300.875

# Row 61
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798

# Row 82
$ws.Range("H82").Value = 742.26666
$ws.Range("I82").Value = 607.8333
$ws.Range("J82").Value = 831.8889
$ws.Range("K82").Value = 607.8333
$ws.Range("L82").Value = 831.8889
$ws.Range("M82").Value = -246.8333
$ws.Range("N82").Value = -1553.8889

# Row 85
$ws.Range("H85").Value = 742.26666
$ws.Range("I85").Value = 607.8333
$ws.Range("J85").Value = 831.8889
$ws.Range("K85").Value = 607.8333
$ws.Range("L85").Value = 831.8889
$ws.Range("M85").Value = 640.1667
$ws.Range("N85").Value = -3327.8889

# Row 100
$ws.Range("H100").Value = 4861
$ws.Range("I100").Value = 4722
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4722
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -4181
$ws.Range("N100").Value = -6082

# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

# Row 122
$ws.Range("H122").Value = 2809.4849
$ws.Range("I122").Value = 2603.3635
$ws.Range("J122").Value = 3221.7273
$ws.Range("K122").Value = 7810.0905
$ws.Range("L122").Value = 9665.1819
$ws.Range("M122").Value = -5360.0905
$ws.Range("N122").Value = -14565.1819

# Row 126
$ws.Range("H126").Value = 6834.407
$ws.Range("I126").Value = 2964.7693
$ws.Range("K126").Value = 8894.3079
$ws.Range("M126").Value = -6424.3079

# Row 132
$ws.Range("H132").Value = 2259.4167
$ws.Range("I132").Value = 1787.9111
$ws.Range("J132").Value = 3673.9333
$ws.Range("K132").Value = 5363.7333
$ws.Range("L132").Value = 11021.7999
$ws.Range("M132").Value = -2833.7333
$ws.Range("N132").Value = -16081.7999


$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2819.45
$ws.Range("I126").Value = 2471.8333
$ws.Range("J126").Value = 3340.875
$ws.Range("K126").Value = 7415.499899999999
$ws.Range("L126").Value = 10022.625
$ws.Range("M126").Value = -4945.499899999999
$ws.Range("N126").Value = -14962.625

# Row 132
$ws.Range("H132").Value = 1897.1224
$ws.Range("I132").Value = 1349.3636
$ws.Range("J132").Value = 3026.875
$ws.Range("K132").Value = 4048.0908
$ws.Range("L132").Value = 9080.625
$ws.Range("M132").Value = -1518.0908
$ws.Range("N132").Value = -14140.625

# Row 136
$ws.Range("H136").Value = 44004.39
$ws.Range("I136").Value = 45854.59
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 137563.77
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -135013.77
$ws.Range("N136").Value = -15000

